$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 1000 (this shifts existing rows 1000..1049 down to 1001..1050)
$ws.Rows.Item(1000).Insert()

# Populate the newly inserted row 1000 with the new record's data
$ws.Cells.Item(1000, 1).Value = 4
$ws.Cells.Item(1000, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(1000, 3).Value = "Los Lagos"
$ws.Cells.Item(1000, 4).Value = 45267
$ws.Cells.Item(1000, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1000, 5).Value = 10
$ws.Cells.Item(1000, 6).Value = 100112033
$ws.Cells.Item(1000, 7).Value = "Lechuga"
$ws.Cells.Item(1000, 8).Value = "Escarola"
$ws.Cells.Item(1000, 9).Value = "Primera"
$ws.Cells.Item(1000, 10).Value = 300
$ws.Cells.Item(1000, 11).Value = 16000
$ws.Cells.Item(1000, 12).Value = 16000
$ws.Cells.Item(1000, 13).Value = 16000
$ws.Cells.Item(1000, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(1000, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(1000, 16).Value = 1067
$ws.Cells.Item(1000, 17).Value = 15
$ws.Cells.Item(1000, 18).Value = "Hortaliza"
